# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect the latest scrape, per commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F5").Value  = 1830
$wsExpo.Range("F6").Value  = 468
$wsExpo.Range("F8").Value  = 165
$wsExpo.Range("F9").Value  = 2389
$wsExpo.Range("F10").Value = 134
$wsExpo.Range("F11").Value = 74
$wsExpo.Range("F13").Value = 1451
$wsExpo.Range("F14").Value = 509
$wsExpo.Range("F20").Value = 198
$wsExpo.Range("F24").Value = 101
$wsExpo.Range("F25").Value = 37
$wsExpo.Range("F26").Value = 1498
$wsExpo.Range("F28").Value = 374
$wsExpo.Range("F29").Value = 277

# --- Sheet "全部类型" (all types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value  = 1830
$wsAll.Range("F7").Value  = 468
$wsAll.Range("F9").Value  = 165
$wsAll.Range("F10").Value = 2389
$wsAll.Range("F11").Value = 134
$wsAll.Range("F12").Value = 74
$wsAll.Range("F14").Value = 1451
$wsAll.Range("F15").Value = 509
$wsAll.Range("F21").Value = 198
$wsAll.Range("F25").Value = 101
$wsAll.Range("F26").Value = 37
$wsAll.Range("F27").Value = 1499
$wsAll.Range("F29").Value = 374
$wsAll.Range("F30").Value = 277
